# "ajout de la colonne ordre"
# Adds a new "Ordre" column (AH) to the "Exigences" sheet, mirroring the
# header/data-row formatting already used by the preceding column (AG),
# widens the new column to match, and updates the selection to reflect
# where the user ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exigences")

# --- Header cell (row 1): copy AG1's formatting onto AH1, then set its text.
$ws.Range("AG1").Copy() | Out-Null
$ws.Range("AH1").PasteSpecial(-4122) | Out-Null
$ws.Range("AH1").Value = "Ordre"

# --- Data cell (row 2): copy AG2's formatting onto AH2, then set its value.
$ws.Range("AG2").Copy() | Out-Null
$ws.Range("AH2").PasteSpecial(-4122) | Out-Null
$ws.Range("AH2").Value = 1

$excel.CutCopyMode = 0

# --- Give the new column roughly the same width as its neighbour.
$ws.Columns.Item(34).ColumnWidth = 13

# --- Reflect the post-edit selection/scroll position.
$ws.Range("AH7").Select() | Out-Null
